$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format so the numeric-looking strings
# are stored as text (matching the inline/shared string cell type used
# by the source data) rather than being coerced into numeric cells.
$ws.Range("B4:E13").NumberFormat = "@"

$ws.Range("B4").Value = "10664.88"
$ws.Range("C4").Value = "9411.33"
$ws.Range("D4").Value = "18822.65"
$ws.Range("E4").Value = "27528.19"

$ws.Range("B5").Value = "10241.10"
$ws.Range("C5").Value = "9021.84"
$ws.Range("D5").Value = "18043.68"
$ws.Range("E5").Value = "26388.89"

$ws.Range("B6").Value = "12289.31"
$ws.Range("C6").Value = "10826.18"
$ws.Range("D6").Value = "21652.36"
$ws.Range("E6").Value = "31666.62"

$ws.Range("B7").Value = "13852.50"
$ws.Range("C7").Value = "11852.31"
$ws.Range("D7").Value = "23704.60"
$ws.Range("E7").Value = "34667.95"

$ws.Range("B8").Value = "15557.11"
$ws.Range("C8").Value = "17354.69"
$ws.Range("D8").Value = "34709.42"
$ws.Range("E8").Value = "50762.45"

$ws.Range("B9").Value = "12658.57"
$ws.Range("C9").Value = "11104.74"
$ws.Range("D9").Value = "22209.48"
$ws.Range("E9").Value = "32481.39"

$ws.Range("B10").Value = "10214.77"
$ws.Range("C10").Value = "8968.28"
$ws.Range("D10").Value = "17936.53"
$ws.Range("E10").Value = "26232.15"

$ws.Range("B11").Value = "14269.39"
$ws.Range("C11").Value = "12492.01"
$ws.Range("D11").Value = "24984.01"
$ws.Range("E11").Value = "36539.14"

$ws.Range("B12").Value = "18236.54"
$ws.Range("C12").Value = "15481.71"
$ws.Range("D12").Value = "30963.40"
$ws.Range("E12").Value = "45283.91"

$ws.Range("B13").Value = "22277.03"
$ws.Range("C13").Value = "18896.72"
$ws.Range("D13").Value = "37793.43"
$ws.Range("E13").Value = "55272.86"

# Drop the temporary Text number format so the cells end up with no
# explicit style (matching the original, unstyled data cells) while
# keeping the values stored as text.
$ws.Range("B4:E13").ClearFormats()
